# Generate Report for Handback
# Updates the "Latest Target File" / "Latest Handback File" / "Latest Handback DateTime"
# and "Error Detail" columns for the zh-cn and de-de handback rows, reflecting that the
# handback for this file came back with an older version than what was already handed
# back (hash 2e769d2766f9f27b7f6478a5ccd03f2dcd4358fb).

$wb = $excel.ActiveWorkbook

$errorDetail = "The handback version of file with file hash 2e769d2766f9f27b7f6478a5ccd03f2dcd4358fb is lower than latest handed back file."
$staleDate = "0001-01-01 00:00:00"
$sourceFileName = "0e5732f1-4137-4e33-9bb7-8d05726bb25f.md"
$sourceFileUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/889435f84ce36066329fa101dbec88d056ac40ca/e2e/0e5732f1-4137-4e33-9bb7-8d05726bb25f.md"

function Clear-HandbackRow($sheetName) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Latest Target File / Latest Handback File are no longer valid - clear them.
    $ws.Range("I2").Value = ""
    $ws.Range("J2").Value = ""

    # Latest Handback DateTime resets to the zero-date sentinel.
    $ws.Range("K2").Value = $staleDate

    # Surface the handback error in the report.
    $ws.Range("P2").Value = $errorDetail

    # Drop the now-stale hyperlink on I2 (the "Latest Target File" link). The sheet's
    # only other hyperlink is A2's "Source File Name" link, which is recreated as-is
    # since the collection doesn't support removing a single member in place.
    $ws.Hyperlinks.Delete()
    $ws.Hyperlinks.Add($ws.Range("A2"), $sourceFileUrl, "", "", $sourceFileName)

    # I2 no longer carries a hyperlink, so drop the HyperLink cell style.
    $ws.Range("I2").Style = "Normal"

    # Column widths: I/J shrink back from the "link text" width now that they're
    # empty, P grows to fit the long error message.
    $ws.Columns.Item(9).ColumnWidth = 17.833333333333332
    $ws.Columns.Item(10).ColumnWidth = 20.833333333333332
    $ws.Columns.Item(16).ColumnWidth = 39.166666666666664
}

Clear-HandbackRow "zh-cn"
Clear-HandbackRow "de-de"
